$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: a new "Unnamed: 0.4" header column is inserted before the old
# "Unnamed: 0.3" column, pushing the existing "Unnamed: 0.x" labels one
# column to the right. What used to be the numeric F1 (0) becomes the text
# label "Unnamed: 0".
$ws.Range("B1").Value = "Unnamed: 0.4"
$ws.Range("C1").Value = "Unnamed: 0.3"
$ws.Range("D1").Value = "Unnamed: 0.2"
$ws.Range("E1").Value = "Unnamed: 0.1"
$ws.Range("F1").Value = "Unnamed: 0"

# --- Row 2 updates
$ws.Range("A2").Value = 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("F2").ClearContents()

# --- Row 3 updates
$ws.Range("A3").Value = 1
$ws.Range("E3").Value = 0

# --- Row 4: a new transaction row, "Pay from Work"
$ws.Range("A4").Value = 3
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Pay from Work"
$ws.Range("H4").Value = 1500
$ws.Range("I4").Value = 0

# J4 needs to hold the literal text "1/14/2025" rather than being
# auto-converted into a date serial number, matching how J2/J3 already
# store their dates as plain text.
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "1/14/2025"
$ws.Range("J4").Style = "Normal"

# Give the new A4 index cell the same bold / thin-border / centered-top
# look as the other index cells in column A (A1:E1, A2, A3).
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").Borders.Weight = 2
